$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "64.453.57"
$ws.Cells.Item(2, 5).Value = "  -0.46%  "

$ws.Cells.Item(3, 4).Value = "3.157.26"
$ws.Cells.Item(3, 5).Value = "  -0.48%  "

$ws.Cells.Item(4, 5).Value = "  +0.11%  "

$ws.Cells.Item(5, 4).Value = "'612.39"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.23%  "

$ws.Cells.Item(6, 4).Value = "'148.20"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -2.16%  "

$ws.Cells.Item(7, 5).Value = "  -0.06%  "

$ws.Cells.Item(8, 4).Value = "3.153.93"
$ws.Cells.Item(8, 5).Value = "  -0.48%  "

$ws.Cells.Item(9, 4).Value = "'0.526"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.76%  "

$ws.Cells.Item(10, 5).Value = "  -0.71%  "

$ws.Cells.Item(11, 4).Value = "'5.47"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -2.12%  "

$ws.Cells.Item(12, 5).Value = "  -0.39%  "

$ws.Cells.Item(13, 4).Value = "'0.0000259"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.06%  "

$ws.Cells.Item(14, 4).Value = "'35.61"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -3.35%  "

$ws.Cells.Item(15, 4).Value = "3.677.96"
$ws.Cells.Item(15, 5).Value = "  -0.41%  "

$ws.Cells.Item(16, 5).Value = "  +2.80%  "

$ws.Cells.Item(17, 4).Value = "64.381.35"
$ws.Cells.Item(17, 5).Value = "  -0.55%  "

$ws.Cells.Item(18, 4).Value = "3.158.18"
$ws.Cells.Item(18, 5).Value = "  -0.44%  "

$ws.Cells.Item(19, 4).Value = "'6.91"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -1.58%  "

$ws.Cells.Item(20, 4).Value = "'479.09"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.20%  "

$ws.Cells.Item(21, 4).Value = "'14.70"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.57%  "

$ws.Cells.Item(22, 4).Value = "'0.716"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.36%  "

$ws.Cells.Item(23, 5).Value = "  +4.11%  "

$ws.Cells.Item(24, 4).Value = "'13.72"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.74%  "

$ws.Cells.Item(25, 5).Value = "  -0.09%  "

$ws.Cells.Item(26, 5).Value = "  -0.02%  "

$ws.Cells.Item(27, 4).Value = "'2.85"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -2.23%  "

$ws.Cells.Item(28, 4).Value = "'8.61"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.53%  "

$ws.Cells.Item(29, 5).Value = "  +0.92%  "

$ws.Cells.Item(30, 5).Value = "  -0.29%  "

$ws.Cells.Item(31, 5).Value = "  -6.49%  "

$ws.Cells.Item(32, 5).Value = "  +0.25%  "

$ws.Cells.Item(33, 5).Value = "  -0.07%  "

$ws.Cells.Item(34, 4).Value = "'26.37"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -1.79%  "

$ws.Cells.Item(35, 5).Value = "  +2.65%  "

$ws.Cells.Item(36, 4).Value = "0.0₃0800"
$ws.Cells.Item(36, 5).Value = "  +8.51%  "

$ws.Cells.Item(37, 5).Value = "  -1.58%  "

$ws.Cells.Item(38, 4).Value = "'3.26"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +2.30%  "

$ws.Cells.Item(39, 4).Value = "'53.19"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -2.69%  "

$ws.Cells.Item(40, 4).Value = "'463.66"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +1.10%  "

$ws.Cells.Item(41, 5).Value = "  -0.36%  "

$ws.Cells.Item(42, 4).Value = "'0.121"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -3.57%  "

$ws.Cells.Item(43, 5).Value = "  -1.08%  "

$ws.Cells.Item(44, 4).Value = "2.853.54"
$ws.Cells.Item(44, 5).Value = "  -0.78%  "

$ws.Cells.Item(45, 4).Value = "'2.33"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -3.39%  "

$ws.Cells.Item(46, 5).Value = "  -1.97%  "

$ws.Cells.Item(47, 5).Value = "  +6.09%  "

$ws.Cells.Item(48, 4).Value = "'26.63"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -1.59%  "

$ws.Cells.Item(49, 4).Value = "'0.998"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.09%  "

$ws.Cells.Item(50, 5).Value = "  -1.50%  "

$ws.Cells.Item(51, 4).Value = "'120.15"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.21%  "
